$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 541.5625
$ws.Range("I5").Value = 355.83334
$ws.Range("J5").Value = 1098.75
$ws.Range("K5").Value = 355.83334
$ws.Range("L5").Value = 1098.75
$ws.Range("M5").Value = -240.83334
$ws.Range("N5").Value = -1328.75
$ws.Range("H6").Value = 858.125
$ws.Range("I6").Value = 1047.5
$ws.Range("K6").Value = 3142.5
$ws.Range("M6").Value = -3030.5
$ws.Range("H92").Value = 867.8823
$ws.Range("I92").Value = 943.3333
$ws.Range("J92").Value = 686.8
$ws.Range("K92").Value = 943.3333
$ws.Range("L92").Value = 686.8
$ws.Range("M92").Value = 304.6667
$ws.Range("N92").Value = -3182.8
$ws.Range("H96").Value = 357.46155
$ws.Range("I96").Value = 222.9
$ws.Range("J96").Value = 806
$ws.Range("K96").Value = 668.7
$ws.Range("L96").Value = 2418
$ws.Range("M96").Value = 704.3
$ws.Range("N96").Value = -5164
$ws.Range("H100").Value = 4065.7144
$ws.Range("J100").Value = 1250
$ws.Range("L100").Value = 1250
$ws.Range("N100").Value = -2332
$ws.Range("H116").Value = 5449.5
$ws.Range("I116").Value = 5449.5
$ws.Range("K116").Value = 5449.5
$ws.Range("M116").Value = -2007.5
$ws.Range("H125").Value = 1000
$ws.Range("J125").Value = 1500
$ws.Range("L125").Value = 13500
$ws.Range("N125").Value = -18420

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 54.25
$ws.Range("I4").Value = 54.25
$ws.Range("K4").Value = 54.25
$ws.Range("M4").Value = 61.75
$ws.Range("H32").Value = 10015.134
$ws.Range("I32").Value = 10015.134
$ws.Range("K32").Value = 10015.134
$ws.Range("M32").Value = -9728.134
$ws.Range("H74").Value = 7083.6665
$ws.Range("I74").Value = 6818.636
$ws.Range("K74").Value = 6818.636
$ws.Range("M74").Value = -5944.636
$ws.Range("H77").Value = 7083.6665
$ws.Range("I77").Value = 6818.636
$ws.Range("K77").Value = 34093.18
$ws.Range("M77").Value = -29725.18
$ws.Range("H97").Value = 374.7143
$ws.Range("I97").Value = 397
$ws.Range("J97").Value = 319
$ws.Range("K97").Value = 397
$ws.Range("L97").Value = 319
$ws.Range("M97").Value = 99
$ws.Range("N97").Value = -1311
$ws.Range("H110").Value = 2469.0454
$ws.Range("I110").Value = 1611.3125
$ws.Range("J110").Value = 4756.3335
$ws.Range("K110").Value = 1611.3125
$ws.Range("L110").Value = 4756.3335
$ws.Range("M110").Value = 433.6875
$ws.Range("N110").Value = -8846.333500000001
$ws.Range("H122").Value = 2284.6
$ws.Range("I122").Value = 2397.4546
$ws.Range("K122").Value = 7192.3638
$ws.Range("M122").Value = -4742.3638

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 30000
$ws.Range("J16").Value = 30000
$ws.Range("L16").Value = 30000
$ws.Range("N16").Value = -30340
$ws.Range("H94").Value = 2792.2727
$ws.Range("I94").Value = 2145.2856
$ws.Range("J94").Value = 3924.5
$ws.Range("K94").Value = 2145.2856
$ws.Range("L94").Value = 3924.5
$ws.Range("M94").Value = -1694.2856
$ws.Range("N94").Value = -4826.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2132
$ws.Range("I58").Value = 2398.5
$ws.Range("J58").Value = 1599
$ws.Range("K58").Value = 2398.5
$ws.Range("L58").Value = 1599
$ws.Range("M58").Value = -2195.5
$ws.Range("N58").Value = -2005
$ws.Range("H86").Value = 6481.6665
$ws.Range("I86").Value = 8247.75
$ws.Range("J86").Value = 2949.5
$ws.Range("K86").Value = 8247.75
$ws.Range("L86").Value = 2949.5
$ws.Range("M86").Value = -7124.75
$ws.Range("N86").Value = -5195.5
$ws.Range("H89").Value = 6481.6665
$ws.Range("I89").Value = 8247.75
$ws.Range("J89").Value = 2949.5
$ws.Range("K89").Value = 41238.75
$ws.Range("L89").Value = 14747.5
$ws.Range("M89").Value = -35622.75
$ws.Range("N89").Value = -25979.5
$ws.Range("H107").Value = 1185.4117
$ws.Range("I107").Value = 1649.3334
$ws.Range("K107").Value = 1649.3334
$ws.Range("M107").Value = 270.6666
$ws.Range("H132").Value = 2478.1538
$ws.Range("I132").Value = 2246.4443
$ws.Range("J132").Value = 2999.5
$ws.Range("K132").Value = 6739.3329
$ws.Range("L132").Value = 8998.5
$ws.Range("M132").Value = -4209.3329
$ws.Range("N132").Value = -14058.5
$ws.Range("H136").Value = 2132
$ws.Range("I136").Value = 2398.5
$ws.Range("J136").Value = 1599
$ws.Range("K136").Value = 7195.5
$ws.Range("L136").Value = 4797
$ws.Range("M136").Value = -4645.5
$ws.Range("N136").Value = -9897

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1777.6154
$ws.Range("I140").Value = 1777.6154
$ws.Range("K140").Value = 5332.8462
$ws.Range("M140").Value = -152.8462

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 8000
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 8000
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 8000
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -8970

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 2833
$ws.Range("J4").Value = 2833
$ws.Range("L4").Value = 2833
$ws.Range("N4").Value = -3059
$ws.Range("H28").Value = 2833
$ws.Range("J28").Value = 2833
$ws.Range("L28").Value = 2833
$ws.Range("N28").Value = -3297
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("N29").ClearContents()
$ws.Range("H37").Value = 2833
$ws.Range("J37").Value = 2833
$ws.Range("L37").Value = 2833
$ws.Range("N37").Value = -3047
$ws.Range("H40").Value = 3339.8
$ws.Range("I40").Value = 3339.8
$ws.Range("K40").Value = 3339.8
$ws.Range("M40").Value = -3203.8
$ws.Range("H46").Value = 3108.75
$ws.Range("I46").Value = 2328.3333
$ws.Range("K46").Value = 2328.3333
$ws.Range("M46").Value = -2140.3333
$ws.Range("H93").Value = 4251.25
$ws.Range("J93").Value = 4001
$ws.Range("L93").Value = 4001
$ws.Range("N93").Value = -6497
$ws.Range("H122").Value = 4332.1665
$ws.Range("I122").Value = 4332.1665
$ws.Range("K122").Value = 12996.4995
$ws.Range("M122").Value = -10546.4995
$ws.Range("H136").Value = 2598.2
$ws.Range("I136").Value = 2598.2
$ws.Range("K136").Value = 7794.599999999999
$ws.Range("M136").Value = -5244.599999999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3536.8
$ws.Range("I96").Value = 3793.111
$ws.Range("J96").Value = 1230
$ws.Range("K96").Value = 3793.111
$ws.Range("L96").Value = 1230
$ws.Range("M96").Value = -2420.111
$ws.Range("N96").Value = -3976
$ws.Range("H107").Value = 979.4
$ws.Range("I107").Value = 778.2222
$ws.Range("J107").Value = 1144
$ws.Range("K107").Value = 2334.6666
$ws.Range("L107").Value = 3432
$ws.Range("M107").Value = -414.6666
$ws.Range("N107").Value = -7272
$ws.Range("H126").Value = 1292.5172
$ws.Range("J126").Value = 1608.0834
$ws.Range("L126").Value = 4824.2502
$ws.Range("N126").Value = -9764.2502
